# Automatic update of files.
#
# 1) Every "Förändrad" (C column) date in the data rows (2..199) moves from
#    2023-09-23 (serial 45192) to 2023-10-03 (serial 45202).
# 2) A new finding ("Knärot") is recorded against case "A 34759-2023", which
#    bumps it above "A 17626-2023" in the sheet's ordering: row 6 now holds
#    the (updated) "A 34759-2023" record and row 7 holds the "A 17626-2023"
#    record that used to occupy row 6 (all other rows keep their row number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1) bump every "Förändrad" date from 45192 to 45202 --------------------
$ws.Range("C2:C199").Value = 45202

# ---- 2) rebuild rows 6 & 7 --------------------------------------------------

# Row 6: "A 34759-2023" (moved up from row 7, with new data merged in)
$ws.Cells.Item(6, 1).Value = "A 34759-2023"
$ws.Cells.Item(6, 2).Value = 45139
$ws.Cells.Item(6, 3).Value = 45202
$ws.Cells.Item(6, 4).Value = "VÄSTMANLANDS LÄN"
$ws.Cells.Item(6, 5).Value = "SURAHAMMAR"
$ws.Cells.Item(6, 6).Value = "Bergvik skog väst AB"
$ws.Cells.Item(6, 7).Value = 38.3
$ws.Cells.Item(6, 8).Value = 5
$ws.Cells.Item(6, 9).Value = 5
$ws.Cells.Item(6, 10).Value = 5
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).Value = 0
$ws.Cells.Item(6, 14).Value = 0
$ws.Cells.Item(6, 15).Value = 6
$ws.Cells.Item(6, 16).Value = 1
$ws.Cells.Item(6, 17).Value = 14
$ws.Cells.Item(6, 18).Value = "Knärot`r`nBlå taggsvamp`r`nGrantaggsvamp`r`nMotaggsvamp`r`nSkogshare`r`nUllticka`r`nDropptaggsvamp`r`nGrönpyrola`r`nMindre märgborre`r`nPlattlummer`r`nVedticka`r`nLopplummer`r`nMattlummer`r`nRevlummer"

$ws.Cells.Item(6, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/artfynd/A 34759-2023.xlsx", "A 34759-2023")'
$ws.Cells.Item(6, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/kartor/A 34759-2023.png", "A 34759-2023")'
$ws.Cells.Item(6, 21).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/knärot/A 34759-2023.png", "A 34759-2023")'
$ws.Cells.Item(6, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/klagomål/A 34759-2023.docx", "A 34759-2023")'
$ws.Cells.Item(6, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/klagomålsmail/A 34759-2023.docx", "A 34759-2023")'
$ws.Cells.Item(6, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/tillsyn/A 34759-2023.docx", "A 34759-2023")'
$ws.Cells.Item(6, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/tillsynsmail/A 34759-2023.docx", "A 34759-2023")'

# Row 7: "A 17626-2023" (moved down from row 6; only the date changed)
$ws.Cells.Item(7, 1).Value = "A 17626-2023"
$ws.Cells.Item(7, 2).Value = 45036
$ws.Cells.Item(7, 3).Value = 45202
$ws.Cells.Item(7, 4).Value = "VÄSTMANLANDS LÄN"
$ws.Cells.Item(7, 5).Value = "SURAHAMMAR"
$ws.Cells.Item(7, 6).ClearContents()
$ws.Cells.Item(7, 7).Value = 5.6
$ws.Cells.Item(7, 8).Value = 3
$ws.Cells.Item(7, 9).Value = 7
$ws.Cells.Item(7, 10).Value = 3
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = 0
$ws.Cells.Item(7, 14).Value = 0
$ws.Cells.Item(7, 15).Value = 4
$ws.Cells.Item(7, 16).Value = 1
$ws.Cells.Item(7, 17).Value = 13
$ws.Cells.Item(7, 18).Value = "Knärot`r`nGranticka`r`nGropticka`r`nUllticka`r`nBlomkålssvamp`r`nGrönpyrola`r`nGullgröppa`r`nKamjordstjärna`r`nStubbspretmossa`r`nTrådticka`r`nVedticka`r`nMattlummer`r`nRevlummer"

$ws.Cells.Item(7, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/artfynd/A 17626-2023.xlsx", "A 17626-2023")'
$ws.Cells.Item(7, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/kartor/A 17626-2023.png", "A 17626-2023")'
$ws.Cells.Item(7, 21).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/knärot/A 17626-2023.png", "A 17626-2023")'
$ws.Cells.Item(7, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/klagomål/A 17626-2023.docx", "A 17626-2023")'
$ws.Cells.Item(7, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/klagomålsmail/A 17626-2023.docx", "A 17626-2023")'
$ws.Cells.Item(7, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/tillsyn/A 17626-2023.docx", "A 17626-2023")'
$ws.Cells.Item(7, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/tillsynsmail/A 17626-2023.docx", "A 17626-2023")'

$wb.Save()
